$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info (row 2-3)
$ws.Range("C2").Value = "Hartmut"

# B3 holds a numeric-looking card number that must remain text;
# format the cell as Text before assigning so it isn't coerced to a number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line (row 5)
$ws.Range("D5").Value = "KONTOSTAND AM 29.04.2025"

# Row 6
$ws.Range("B6").Value = "01.05."
$ws.Range("C6").Value = "02.05."
$ws.Range("D6").Value = "PAYPAL FZLZYR"
$ws.Range("E6").Value = "91,25-"

# Row 7
$ws.Range("B7").Value = "02.05."
$ws.Range("C7").Value = "03.05."
$ws.Range("D7").Value = "PAYPAL PGPPWJ"
$ws.Range("E7").Value = "15,65-"

# Row 8
$ws.Range("B8").Value = "03.05."
$ws.Range("C8").Value = "04.05."
$ws.Range("D8").Value = "KARTENZ./03.05 REWE RO"
$ws.Range("E8").Value = "146,21-"

# Row 9
$ws.Range("B9").Value = "05.05."
$ws.Range("C9").Value = "06.05."
$ws.Range("D9").Value = "EBAY MKTPLC EU XSFESG"
$ws.Range("E9").Value = "78,67-"

# Row 10 - previously blank trailing row, now filled with another transaction
$ws.Range("B10").Value = "07.05."
$ws.Range("C10").Value = "08.05."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-57788630"
$ws.Range("E10").Value = "52,72-"
# Match the right-aligned, non-wrapping style used by the other amount cells
# in this column (E6:E9) instead of the blank-row style it inherited.
$ws.Range("E10").HorizontalAlignment = -4152  # xlRight
$ws.Range("E10").VerticalAlignment = -4107    # xlBottom
$ws.Range("E10").WrapText = $false

# Closing balance line (row 12)
$ws.Range("D12").Value = "KONTOSTAND AM 12.05.2025"
$ws.Range("E12").Value = "384,50-"

# Next statement date (row 13)
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 22.05.2025"
